# feat: add 2022-Q3 data
#
# 1. "总计" sheet: insert a new summary row for 2022-Q3 above the existing
#    2022-Q2 row (so 2022-Q3 is first, 2022-Q2 shifts down to row 3).
# 2. Add a new "2022-Q3" worksheet (positioned before "2022-Q2", after
#    "总计") holding the per-fund breakdown for the quarter, built by
#    duplicating the existing "2022-Q2" sheet (so header styling / column
#    layout matches exactly) and then overwriting it with the new data.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("总计")
$ws2 = $wb.Worksheets.Item("2022-Q2")

# ---------------------------------------------------------------------
# Step 1: "总计" sheet - insert the 2022-Q3 summary row
# ---------------------------------------------------------------------

# Remember the current (2022-Q2) row-2 values before they get pushed down.
$oldDate  = $ws1.Range("B2").Value2
$oldCount = $ws1.Range("C2").Value2
$oldValue = $ws1.Range("D2").Value2

# Push the existing data row down to row 3.
$ws1.Rows.Item(2).Insert()
$ws1.Range("B2:D2").ClearFormats()

# Restore the (now row 3) 2022-Q2 data, with its running index updated.
$ws1.Range("A3").Value = 1
$ws1.Range("B3").Value = $oldDate
$ws1.Range("C3").Value = $oldCount
$ws1.Range("D3").Value = $oldValue
$ws1.Range("B1").Copy()
$ws1.Range("A3").PasteSpecial(-4122)

# Fill in the new 2022-Q3 summary row.
$ws1.Range("A2").Value = 0
$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 8
$ws1.Range("D2").Value = 0.74
$ws1.Range("B1").Copy()
$ws1.Range("A2").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 2: add the "2022-Q3" worksheet (fund-level detail)
# ---------------------------------------------------------------------

# Duplicate "2022-Q2" (placed immediately before it) so the new sheet
# starts out with identical column headers / styles, then rename it.
$ws2.Copy($ws2)
$ws3 = $wb.Worksheets.Item("2022-Q2 (2)")
$ws3.Name = "2022-Q3"

# Fund detail rows for 2022-Q3.
$rows = @(
    @("012719", "华夏新兴经济一年持有混合A", "17.53", "87.21", "2.69", "0.4716", 9),
    @("010994", "博时创新经济混合A",         "3.35",  "91.74", "4.93", "0.1652", 7),
    @("002295", "广发稳安灵活配置混合A",     "1.58",  "69.63", "3.86", "0.0610", 8),
    @("012720", "华夏新兴经济一年持有混合C", "0.98",  "87.21", "2.69", "0.0264", 9),
    @("010995", "博时创新经济混合C",         "0.38",  "91.74", "4.93", "0.0187", 7),
    @("008604", "广发稳安灵活配置混合C",     "0.02",  "69.63", "3.86", "0.0008", 8),
    @("011786", "工银聚安混合A",             "0.90",  "24.17", "0.02", "0.0002", 7),
    @("011787", "工银聚安混合C",             "0.01",  "24.17", "0.02", "0.0000", 7)
)

# Row 2 through 9 all share the same plain (no-border) style as the
# original data row 2 - copy it down so rows 3-9 match row 2's styling
# before writing values into them.
$ws3.Range("A2:H2").Copy()
$ws3.Range("A3:H9").PasteSpecial(-4122)

# Make sure text-like numeric strings (fund code / size / position /
# ratio / market value columns) are stored as text, matching the source
# data (preserves the leading zeros in fund codes, e.g. "012719").
$ws3.Range("B2:B9").NumberFormat = "@"
$ws3.Range("D2:G9").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    $ws3.Range("A$r").Value = $i
    $ws3.Range("B$r").Value = $row[0]
    $ws3.Range("C$r").Value = $row[1]
    $ws3.Range("D$r").Value = $row[2]
    $ws3.Range("E$r").Value = $row[3]
    $ws3.Range("F$r").Value = $row[4]
    $ws3.Range("G$r").Value = $row[5]
    $ws3.Range("H$r").Value = $row[6]
}

# ---------------------------------------------------------------------
# Step 3: restore "总计" as the active sheet (unchanged selection state)
# ---------------------------------------------------------------------
$ws1.Activate()

Write-Host "2022-Q3 data added"
